# Natmi following Dr Hou advice
# Extend the LR-pairs sheet (Tnc -> Ptprz1) from a 2-row (ECs/FAPs-only)
# result to the full 3x3 cluster matrix (ECs, FAPs, sCs), and refresh the
# already-present rows' values to match the new NATMI run. Also introduces
# a new "ECs" cluster-name string used by several rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> Tnc | Ptprz1 -> ECs
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Tnc"
$ws.Range("C2").Value = "Ptprz1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 1.442371333333333
$ws.Range("H2").Value = 4.327114
$ws.Range("I2").Value = 0.03522044016446201
$ws.Range("J2").Value = 0.03522044016446201
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.03995766666666666
$ws.Range("N2").Value = 0.119873
$ws.Range("O2").Value = 0.005314930928687666
$ws.Range("P2").Value = 0.005314930928687667
$ws.Range("Q2").Value = 0.05763379294688888
$ws.Range("R2").Value = 0.518704136522
$ws.Range("S2").Value = 0.0001871942067520925
$ws.Range("T2").Value = 0.0001871942067520925

# Row 3: ECs -> Tnc | Ptprz1 -> FAPs
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Tnc"
$ws.Range("C3").Value = "Ptprz1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 1.442371333333333
$ws.Range("H3").Value = 4.327114
$ws.Range("I3").Value = 0.03522044016446201
$ws.Range("J3").Value = 0.03522044016446201
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.03069133333333333
$ws.Range("N3").Value = 0.092074
$ws.Range("O3").Value = 0.004082378436578614
$ws.Range("P3").Value = 0.004082378436578615
$ws.Range("Q3").Value = 0.04426829938177777
$ws.Range("R3").Value = 0.398414694436
$ws.Range("S3").Value = 0.000143783165454207
$ws.Range("T3").Value = 0.0001437831654542071

# Row 4: ECs -> Tnc | Ptprz1 -> sCs
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Tnc"
$ws.Range("C4").Value = "Ptprz1"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 1.442371333333333
$ws.Range("H4").Value = 4.327114
$ws.Range("I4").Value = 0.03522044016446201
$ws.Range("J4").Value = 0.03522044016446201
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 7.447354000000001
$ws.Range("N4").Value = 22.342062
$ws.Range("O4").Value = 0.9906026906347337
$ws.Range("P4").Value = 0.9906026906347338
$ws.Range("Q4").Value = 10.74184991878533
$ws.Range("R4").Value = 96.67664926906801
$ws.Range("S4").Value = 0.03488946279225572
$ws.Range("T4").Value = 0.03488946279225572

# Row 5: FAPs -> Tnc | Ptprz1 -> ECs
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Tnc"
$ws.Range("C5").Value = "Ptprz1"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 11.331397
$ws.Range("H5").Value = 33.994191
$ws.Range("I5").Value = 0.2766948987373093
$ws.Range("J5").Value = 0.2766948987373092
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.03995766666666666
$ws.Range("N5").Value = 0.119873
$ws.Range("O5").Value = 0.005314930928687666
$ws.Range("P5").Value = 0.005314930928687667
$ws.Range("Q5").Value = 0.4527761841936667
$ws.Range("R5").Value = 4.074985657742999
$ws.Range("S5").Value = 0.001470614275109027
$ws.Range("T5").Value = 0.001470614275109027

# Row 6: FAPs -> Tnc | Ptprz1 -> FAPs
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Tnc"
$ws.Range("C6").Value = "Ptprz1"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 11.331397
$ws.Range("H6").Value = 33.994191
$ws.Range("I6").Value = 0.2766948987373093
$ws.Range("J6").Value = 0.2766948987373092
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.03069133333333333
$ws.Range("N6").Value = 0.092074
$ws.Range("O6").Value = 0.004082378436578614
$ws.Range("P6").Value = 0.004082378436578615
$ws.Range("Q6").Value = 0.3477756824593334
$ws.Range("R6").Value = 3.129981142134
$ws.Range("S6").Value = 0.001129573288116494
$ws.Range("T6").Value = 0.001129573288116495

# Row 7: FAPs -> Tnc | Ptprz1 -> sCs
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Tnc"
$ws.Range("C7").Value = "Ptprz1"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 11.331397
$ws.Range("H7").Value = 33.994191
$ws.Range("I7").Value = 0.2766948987373093
$ws.Range("J7").Value = 0.2766948987373092
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 7.447354000000001
$ws.Range("N7").Value = 22.342062
$ws.Range("O7").Value = 0.9906026906347337
$ws.Range("P7").Value = 0.9906026906347338
$ws.Range("Q7").Value = 84.38892477353801
$ws.Range("R7").Value = 759.5003229618421
$ws.Range("S7").Value = 0.2740947111740838
$ws.Range("T7").Value = 0.2740947111740837

# Row 8: sCs -> Tnc | Ptprz1 -> ECs
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Tnc"
$ws.Range("C8").Value = "Ptprz1"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 28.17890933333333
$ws.Range("H8").Value = 84.536728
$ws.Range("I8").Value = 0.6880846610982287
$ws.Range("J8").Value = 0.6880846610982286
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.03995766666666666
$ws.Range("N8").Value = 0.119873
$ws.Range("O8").Value = 0.005314930928687666
$ws.Range("P8").Value = 0.005314930928687667
$ws.Range("Q8").Value = 1.125963466171555
$ws.Range("R8").Value = 10.133671195544
$ws.Range("S8").Value = 0.003657122446826547
$ws.Range("T8").Value = 0.003657122446826547

# Row 9: sCs -> Tnc | Ptprz1 -> FAPs
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Tnc"
$ws.Range("C9").Value = "Ptprz1"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 28.17890933333333
$ws.Range("H9").Value = 84.536728
$ws.Range("I9").Value = 0.6880846610982287
$ws.Range("J9").Value = 0.6880846610982286
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.03069133333333333
$ws.Range("N9").Value = 0.092074
$ws.Range("O9").Value = 0.004082378436578614
$ws.Range("P9").Value = 0.004082378436578615
$ws.Range("Q9").Value = 0.8648482993191111
$ws.Range("R9").Value = 7.783634693872
$ws.Range("S9").Value = 0.002809021983007913
$ws.Range("T9").Value = 0.002809021983007913

# Row 10: sCs -> Tnc | Ptprz1 -> sCs
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Tnc"
$ws.Range("C10").Value = "Ptprz1"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 28.17890933333333
$ws.Range("H10").Value = 84.536728
$ws.Range("I10").Value = 0.6880846610982287
$ws.Range("J10").Value = 0.6880846610982286
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 7.447354000000001
$ws.Range("N10").Value = 22.342062
$ws.Range("O10").Value = 0.9906026906347337
$ws.Range("P10").Value = 0.9906026906347338
$ws.Range("Q10").Value = 209.8583131392374
$ws.Range("R10").Value = 1888.724818253136
$ws.Range("S10").Value = 0.6816185166683942
$ws.Range("T10").Value = 0.6816185166683942
